$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('leaderboard2')
$ws.Range('D4').Value = 1189
$ws.Range('N9').Value = "'" + '51.'
$ws.Range('O9').Value = 'VelNewt'
$ws.Range('P9').Value = 0
$ws.Range('D12').Value = 946
$ws.Range('B14').Value = 'Dernière update le 18.03.25 à 01:27'

$ws = $wb.Worksheets.Item('leaderboard3')
$ws.Range('F4').Value = '_Linca'
$ws.Range('F5').Value = 'Brybry_'
$ws.Range('O8').Value = 'VelNewt'
$ws.Range('N9').Value = "'" + '51.'
$ws.Range('O9').Value = 'nisqylegoat'
$ws.Range('P9').Value = 0
$ws.Range('L10').Value = 'Horty_'
$ws.Range('L11').Value = 'ARELIANN'
$ws.Range('B14').Value = 'Dernière update le 18.03.25 à 01:27'

$ws = $wb.Worksheets.Item('leaderboard4')
$ws.Range('L3').Value = 'Maxouzboub'
$ws.Range('O6').Value = 'VelNewt'
$ws.Range('O7').Value = 'ZeratoR'
$ws.Range('L8').Value = 'ARELIANN'
$ws.Range('O8').Value = 'Horty_'
$ws.Range('L9').Value = 'Mynth0s'
$ws.Range('N9').Value = "'" + '51.'
$ws.Range('O9').Value = 'nisqylegoat'
$ws.Range('P9').Value = 0
$ws.Range('I13').Value = 'LadySundae'
$ws.Range('B14').Value = 'Dernière update le 18.03.25 à 01:27'

$ws = $wb.Worksheets.Item('leaderboard5')
$ws.Range('L3').Value = 'Terraciid'
$ws.Range('O3').Value = 'HarryLafranc'
$ws.Range('L4').Value = 'ARELIANN'
$ws.Range('L5').Value = 'ZeratoR'
$ws.Range('O5').Value = 'VelNewt'
$ws.Range('I6').Value = 'Onutrem'
$ws.Range('L6').Value = 'XoTrixy'
$ws.Range('O6').Value = 'Horty_'
$ws.Range('I7').Value = 'Kaatsup'
$ws.Range('L7').Value = 'JLTootmy'
$ws.Range('O7').Value = 'TheGuill84'
$ws.Range('L8').Value = 'Pepito_kawazakii'
$ws.Range('O8').Value = 'Grimkujow'
$ws.Range('L9').Value = 'Mynth0s'
$ws.Range('N9').Value = "'" + '51.'
$ws.Range('O9').Value = 'nisqylegoat'
$ws.Range('P9').Value = 0
$ws.Range('I10').Value = 'Elspawn'
$ws.Range('L10').Value = 'Wingobear'
$ws.Range('I11').Value = 'Bytell2'
$ws.Range('L11').Value = 'CrocodyleTV'
$ws.Range('I12').Value = 'Angle_Droit'
$ws.Range('L12').Value = 'BagheraJones'
$ws.Range('L13').Value = 'Hiro_Ammar'
$ws.Range('B14').Value = 'Dernière update le 18.03.25 à 01:27'
